$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cellRef, $value)
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextCell "D2" "41.824.50"
$ws.Range("E2").Value = "  -0.73%  "
Set-TextCell "D3" "2.259.48"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("E4").Value = "  -0.09%  "
Set-TextCell "D5" "304.86"
$ws.Range("E5").Value = "  -0.32%  "
Set-TextCell "D6" "95.15"
$ws.Range("E6").Value = "  +1.79%  "
Set-TextCell "D7" "0.525"
$ws.Range("E7").Value = "  -1.02%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -0.04%  "
Set-TextCell "D10" "34.94"
$ws.Range("E10").Value = "  +6.15%  "
Set-TextCell "D11" "0.0789"
$ws.Range("E11").Value = "  -2.02%  "
$ws.Range("E12").Value = "  -0.55%  "
Set-TextCell "D13" "6.70"
$ws.Range("E13").Value = "  -0.22%  "
Set-TextCell "D14" "2.612.62"
$ws.Range("E14").Value = "  -0.43%  "
$ws.Range("E15").Value = "  +0.13%  "
Set-TextCell "D16" "2.260.58"
$ws.Range("E16").Value = "  -0.53%  "
Set-TextCell "D17" "0.788"
$ws.Range("E17").Value = "  +0.28%  "
Set-TextCell "D18" "41.735.22"
$ws.Range("E18").Value = "  -0.64%  "
Set-TextCell "D19" "12.34"
$ws.Range("E19").Value = "  -3.14%  "
Set-TextCell "D20" "0.0₃0899"
$ws.Range("E20").Value = "  -2.22%  "
$ws.Range("E21").Value = "  -0.79%  "
Set-TextCell "D22" "67.85"
$ws.Range("E22").Value = "  -0.52%  "
Set-TextCell "D23" "236.63"
$ws.Range("E23").Value = "  -3.11%  "
Set-TextCell "D24" "2.55"
$ws.Range("E24").Value = "  -1.89%  "
Set-TextCell "D25" "0.999"
$ws.Range("E25").Value = "  -0.03%  "
Set-TextCell "D26" "1.91"
$ws.Range("E26").Value = "  -1.59%  "
Set-TextCell "D27" "23.58"
$ws.Range("E27").Value = "  -2.09%  "
Set-TextCell "D28" "36.54"
$ws.Range("E28").Value = "  +4.00%  "
$ws.Range("E29").Value = "  +0.56%  "
Set-TextCell "D30" "9.42"
$ws.Range("E30").Value = "  -2.86%  "
Set-TextCell "D31" "159.86"
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextCell "D32" "0.999"
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell "D33" "5.19"
$ws.Range("E33").Value = "  -2.84%  "
Set-TextCell "D34" "3.17"
$ws.Range("E34").Value = "  +4.57%  "
Set-TextCell "D35" "0.0733"
$ws.Range("E35").Value = "  -1.58%  "
Set-TextCell "D36" "16.89"
$ws.Range("E36").Value = "  -1.61%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("E38").Value = "  -1.10%  "
$ws.Range("E39").Value = "  +0.64%  "
$ws.Range("E40").Value = "  -2.52%  "
Set-TextCell "D41" "3.97"
$ws.Range("E41").Value = "  -0.96%  "
$ws.Range("E42").Value = "  +2.49%  "
Set-TextCell "D43" "1.967.86"
$ws.Range("E43").Value = "  -2.30%  "
Set-TextCell "D44" "0.0281"
$ws.Range("E44").Value = "  -0.77%  "
Set-TextCell "D45" "18.76"
$ws.Range("E45").Value = "  -5.00%  "
Set-TextCell "D46" "2.91"
$ws.Range("E46").Value = "  -0.09%  "
Set-TextCell "D47" "9.83"
$ws.Range("E47").Value = "  -4.07%  "
Set-TextCell "D48" "52.67"
$ws.Range("E48").Value = "  -1.24%  "
Set-TextCell "D49" "71.73"
$ws.Range("E49").Value = "  -1.30%  "
$ws.Range("E50").Value = "  -1.56%  "
Set-TextCell "D51" "90.80"
$ws.Range("E51").Value = "  -1.16%  "
